# fix parameterization, run strategy and influence experiments
#
# The influence-score block (rows 22-27, cols C:F) gets shifted up by one
# row: the "highlighted" (red, style 6) -0.5 scores that lived on row 25
# move to row 24, the plain (style 4) -0.5 scores that lived on row 26
# move to row 25, and row 26 is cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 24: take on the "highlighted" style (row 22 already uses it)
#     and the -0.5 values that used to sit on row 25 ---
$ws.Range("C22:F22").Copy()
$ws.Range("C24:F24").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C24:F24").Value = -0.5

# --- Row 25: switch to the plain style (matches B25) and keep -0.5 ---
$ws.Range("B25").Copy()
$ws.Range("C25:F25").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C25:F25").Value = -0.5

# --- Row 26: clear the values that were duplicated onto row 25 ---
$ws.Range("C26:F26").ClearContents()

# --- View state: selection moves to B25 ---
$ws.Range("B25").Select()
